$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Games")

# Append three more rows (63-65) following the same pattern as the
# existing trailing rows (57-62): A = row-1, B = 0, C = 0, D = 3,
# E = "14.01.2020" (stored as text).
for ($row = 63; $row -le 65; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 1
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 3
    $ws.Cells.Item($row, 5).Value = "14.01.2020"
}
